$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates
$ws.Range("D2").Value = 44320
$ws.Range("M2").Value = 50
$ws.Range("N2").Value = 18000
$ws.Range("P2").Value = 18800
$ws.Range("R2").Value = "Provincia de Limarí"
$ws.Range("S2").Value = 1044

# Row 3 updates
$ws.Range("D3").Value = 44362
$ws.Range("M3").Value = 100
$ws.Range("N3").Value = 19000
$ws.Range("P3").Value = 19500
$ws.Range("R3").Value = "Provincia de Curicó"
$ws.Range("S3").Value = 1083
